$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.239.87'
$ws.Range("E2").Value = '  -2.77%  '
$ws.Range("D3").Value = '2.555.09'
$ws.Range("E3").Value = '  -4.52%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Formula = "'544.61"
$ws.Range("E5").Value = '  -0.90%  '
$ws.Range("D6").Formula = "'151.10"
$ws.Range("E6").Value = '  -3.93%  '
$ws.Range("D7").Formula = "'1.00"
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").Formula = "'0.584"
$ws.Range("E8").Value = '  -0.24%  '
$ws.Range("E9").Value = '  -2.80%  '
$ws.Range("E10").Value = '  -0.70%  '
$ws.Range("D11").Formula = "'5.45"
$ws.Range("E11").Value = '  +7.46%  '
$ws.Range("D12").Formula = "'0.360"
$ws.Range("E12").Value = '  -1.96%  '
$ws.Range("D13").Value = '3.008.13'
$ws.Range("E13").Value = '  -4.45%  '
$ws.Range("D14").Formula = "'25.06"
$ws.Range("E14").Value = '  -3.80%  '
$ws.Range("D15").Value = '61.224.02'
$ws.Range("E15").Value = '  -2.60%  '
$ws.Range("E16").Value = '  -1.99%  '
$ws.Range("D17").Value = '2.556.20'
$ws.Range("E17").Value = '  -4.49%  '
$ws.Range("D18").Formula = "'11.44"
$ws.Range("E18").Value = '  -4.01%  '
$ws.Range("D19").Formula = "'4.47"
$ws.Range("E19").Value = '  -1.84%  '
$ws.Range("D20").Formula = "'334.44"
$ws.Range("E20").Value = '  -1.89%  '
$ws.Range("D21").Formula = "'1.00"
$ws.Range("E21").Value = '  +0.38%  '
$ws.Range("D22").Formula = "'5.93"
$ws.Range("E22").Value = '  -5.49%  '
$ws.Range("E23").Value = '  -4.24%  '
$ws.Range("D24").Formula = "'62.54"
$ws.Range("E24").Value = '  -1.19%  '
$ws.Range("D25").Formula = "'0.165"
$ws.Range("E25").Value = '  -1.59%  '
$ws.Range("E26").Value = '  +0.56%  '
$ws.Range("D27").Formula = "'8.00"
$ws.Range("E27").Value = '  -1.18%  '
$ws.Range("D28").Formula = "'7.07"
$ws.Range("E28").Value = '  +1.22%  '
$ws.Range("D29").Value = '0.0₃0812'
$ws.Range("E29").Value = '  -4.27%  '
$ws.Range("E30").Value = '  -2.17%  '
$ws.Range("D31").Formula = "'1.86"
$ws.Range("E31").Value = '  -2.88%  '
$ws.Range("D32").Formula = "'161.72"
$ws.Range("E32").Value = '  -2.32%  '
$ws.Range("E33").Value = '  +0.06%  '
$ws.Range("D34").Formula = "'4.76"
$ws.Range("E34").Value = '  -0.59%  '
$ws.Range("D35").Formula = "'18.90"
$ws.Range("E35").Value = '  -2.84%  '
$ws.Range("E36").Value = '  -3.12%  '
$ws.Range("D37").Formula = "'1.74"
$ws.Range("E37").Value = '  -1.26%  '
$ws.Range("D38").Formula = "'320.91"
$ws.Range("E38").Value = '  -5.49%  '
$ws.Range("D39").Formula = "'5.89"
$ws.Range("E39").Value = '  -4.15%  '
$ws.Range("D40").Formula = "'0.873"
$ws.Range("E40").Value = '  -6.41%  '
$ws.Range("D41").Formula = "'3.85"
$ws.Range("E41").Value = '  -1.63%  '
$ws.Range("D42").Formula = "'37.06"
$ws.Range("E42").Value = '  -2.56%  '
$ws.Range("E43").Value = '  +0.16%  '
$ws.Range("D44").Formula = "'20.29"
$ws.Range("E44").Value = '  -1.90%  '
$ws.Range("D45").Formula = "'10.90"
$ws.Range("E45").Value = '  -1.24%  '
$ws.Range("D46").Formula = "'0.600"
$ws.Range("E46").Value = '  -2.18%  '
$ws.Range("D47").Formula = "'0.0958"
$ws.Range("E47").Value = '  -1.20%  '
$ws.Range("D48").Formula = "'0.0535"
$ws.Range("E48").Value = '  -4.18%  '
$ws.Range("D49").Formula = "'19.15"
$ws.Range("E49").Value = '  -5.23%  '
$ws.Range("D50").Formula = "'0.0234"
$ws.Range("E50").Value = '  -1.87%  '
$ws.Range("D51").Value = '2.031.77'
$ws.Range("E51").Value = '  -2.04%  '
